# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 29 de Abril de 2020 a las 18:52"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1045717
$ws.Range("C4").Value = 9952
$ws.Range("D4").Value = 144048
$ws.Range("E4").Value = 841586
$ws.Range("G4").Value = 817
$ws.Range("H4").Value = 60083

# Row 10 - Turquia
$ws.Range("B10").Value = 117589
$ws.Range("C10").Value = 2936
$ws.Range("D10").Value = 44022
$ws.Range("E10").Value = 70486
$ws.Range("F10").Value = 1574
$ws.Range("G10").Value = 89
$ws.Range("H10").Value = 3081

# Row 15 - Canada
$ws.Range("D15").Value = 19886
$ws.Range("E15").Value = 27583

# Row 48 - Republica Dominicana
$ws.Range("B48").Value = 6652
$ws.Range("C48").Value = 236
$ws.Range("D48").Value = 1228
$ws.Range("E48").Value = 5131
$ws.Range("G48").Value = 7
$ws.Range("H48").Value = 293

# Row 55 - Marruecos
$ws.Range("B55").Value = 4321
$ws.Range("C55").Value = 69
$ws.Range("D55").Value = 928
$ws.Range("E55").Value = 3225
$ws.Range("G55").Value = 3
$ws.Range("H55").Value = 168

# Row 61 - Kazajistan
$ws.Range("B61").Value = 3138
$ws.Range("C61").Value = 111
$ws.Range("D61").Value = 819
$ws.Range("E61").Value = 2294

# Row 109 - Georgia
$ws.Range("D109").Value = 178
$ws.Range("E109").Value = 333

# Row 136 - Birmania
$ws.Range("D136").Value = 27
$ws.Range("E136").Value = 117
